# "new date for fontaine"
# The "Prix de Fontaine" event (previously row 6, dated "Dim 31 Mars") moves up
# to a new row 4 (right after "Dim 3 Mars"), keeping the same date but gaining a
# "Nouvelle date !" note in column F. Rows that used to sit between the old and
# new position (old rows 4-5, "Sam 23 Mars" / "Dim 24 Mars") shift down to make
# room, and the now-duplicate old row is removed so everything below is back in
# its original place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 4, pushing rows 4.. down by one.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the Fontaine event data + the new note.
$ws.Range("A4").Value = "Dim 31 Mars"
$ws.Range("B4").Value = "Prix de Fontaine  "
$ws.Range("C4").Value = "ACT Belfort"
$ws.Range("D4").Value = "Route"
$ws.Range("E4").Value = "fontaine"
$ws.Range("F4").Value = "Nouvelle date !"

# The inserted row copied formatting from the row above it; the target state
# has no explicit style on row 4, so clear it back to Normal.
$ws.Rows.Item(4).ClearFormats()

# The old "Prix de Fontaine" row (originally row 6) is now row 7 after the
# insert above (rows 4/5 shifted to 5/6, old row 6 shifted to 7). Remove it
# since its data now lives in the new row 4.
$ws.Rows.Item(7).Delete()

# Update the selection to match the saved view state.
$ws.Range("A5").Select() | Out-Null
